# Refresh the crypto price/volume table with the latest scraped values
# (GitHub Actions data-refresh commit). Cells that hold a purely numeric
# "Price" string (single decimal point, e.g. "65.63") are force-formatted
# as Text first so Excel doesn't silently coerce them into floating point
# numbers and lose the exact printed representation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.944.37"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "1.672.20"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("E6").Value = "  +1.45%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E9").Value = "  +0.39%  "
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0889"
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("D12").Value = "1.907.54"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("D13").Value = "1.636.59"
$ws.Range("E13").Value = "  -1.15%  "
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.63"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").Value = "26.962.40"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.08"
$ws.Range("E18").Value = "  +3.94%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "234.47"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("D20").Value = "0.0₃0733"
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("E23").Value = "  -1.44%  "
$ws.Range("E24").Value = "  -1.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.64"
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.15"
$ws.Range("E26").Value = "  +0.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.97"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("E28").Value = "  -1.49%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("E32").Value = "  +0.91%  "
$ws.Range("D33").Value = "1.471.09"
$ws.Range("E33").Value = "  -5.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.14"
$ws.Range("E34").Value = "  +2.31%  "
$ws.Range("E35").Value = "  +1.92%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("E38").Value = "  -0.80%  "
$ws.Range("E39").Value = "  +1.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.06"
$ws.Range("E40").Value = "  +7.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.80"
$ws.Range("E41").Value = "  -3.81%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("E43").Value = "  +2.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "66.67"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").Value = "1.815.31"
$ws.Range("E45").Value = "  +1.03%  "
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.34"
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.71"
$ws.Range("E51").Value = "  +0.42%  "
